{"js": "// Update the date line and all 100 two-digit multiplication problems in the\n// practice-sheet table. Every `w:t` value in this document is unique, so each\n// [oldText, newText] pair below identifies exactly one run; we search for\n// oldText and replace its text with newText, once per pair, in document order.\nconst replacements = [\n  [\"2023-04-26 Wednesday\", \"2023-04-27 Thursday\"],\n  [\"96\u00d763=\", \"27\u00d755=\"],\n  [\"88\u00d716=\", \"86\u00d789=\"],\n  [\"46\u00d777=\", \"91\u00d744=\"],\n  [\"87\u00d748=\", \"93\u00d742=\"],\n  [\"50\u00d740=\", \"91\u00d757=\"],\n  [\"84\u00d791=\", \"85\u00d797=\"],\n  [\"48\u00d740=\", \"43\u00d7100=\"],\n  [\"84\u00d745=\", \"12\u00d796=\"],\n  [\"27\u00d794=\", \"22\u00d785=\"],\n  [\"94\u00d747=\", \"25\u00d767=\"],\n  [\"43\u00d777=\", \"99\u00d737=\"],\n  [\"11\u00d745=\", \"48\u00d762=\"],\n  [\"55\u00d751=\", \"87\u00d766=\"],\n  [\"54\u00d717=\", \"71\u00d759=\"],\n  [\"10\u00d785=\", \"11\u00d799=\"],\n  [\"61\u00d796=\", \"10\u00d723=\"],\n  [\"50\u00d754=\", \"30\u00d790=\"],\n  [\"98\u00d711=\", \"60\u00d779=\"],\n  [\"62\u00d746=\", \"55\u00d764=\"],\n  [\"10\u00d726=\", \"74\u00d784=\"],\n  [\"10\u00d729=\", \"39\u00d726=\"],\n  [\"51\u00d745=\", \"39\u00d733=\"],\n  [\"62\u00d733=\", \"33\u00d782=\"],\n  [\"29\u00d733=\", \"57\u00d750=\"],\n  [\"27\u00d722=\", \"13\u00d799=\"],\n  [\"56\u00d758=\", \"48\u00d730=\"],\n  [\"41\u00d760=\", \"11\u00d763=\"],\n  [\"88\u00d736=\", \"83\u00d763=\"],\n  [\"16\u00d753=\", \"20\u00d755=\"],\n  [\"52\u00d767=\", \"96\u00d711=\"],\n  [\"54\u00d764=\", \"94\u00d791=\"],\n  [\"35\u00d771=\", \"74\u00d771=\"],\n  [\"87\u00d783=\", \"54\u00d773=\"],\n  [\"39\u00d720=\", \"49\u00d794=\"],\n  [\"19\u00d748=\", \"36\u00d779=\"],\n  [\"12\u00d749=\", \"47\u00d768=\"],\n  [\"82\u00d758=\", \"97\u00d750=\"],\n  [\"11\u00d754=\", \"74\u00d775=\"],\n  [\"70\u00d711=\", \"59\u00d738=\"],\n  [\"28\u00d714=\", \"71\u00d763=\"],\n  [\"73\u00d734=\", \"65\u00d771=\"],\n  [\"78\u00d762=\", \"52\u00d715=\"],\n  [\"32\u00d745=\", \"73\u00d712=\"],\n  [\"85\u00d727=\", \"80\u00d743=\"],\n  [\"81\u00d793=\", \"95\u00d719=\"],\n  [\"89\u00d740=\", \"96\u00d748=\"],\n  [\"53\u00d744=\", \"12\u00d718=\"],\n  [\"77\u00d723=\", \"20\u00d732=\"],\n  [\"77\u00d715=\", \"42\u00d736=\"],\n  [\"78\u00d782=\", \"61\u00d791=\"],\n  [\"84\u00d769=\", \"38\u00d778=\"],\n  [\"38\u00d796=\", \"39\u00d731=\"],\n  [\"44\u00d767=\", \"60\u00d758=\"],\n  [\"43\u00d788=\", \"67\u00d756=\"],\n  [\"83\u00d798=\", \"37\u00d795=\"],\n  [\"98\u00d784=\", \"64\u00d777=\"],\n  [\"86\u00d774=\", \"71\u00d733=\"],\n  [\"92\u00d794=\", \"34\u00d766=\"],\n  [\"48\u00d720=\", \"30\u00d735=\"],\n  [\"17\u00d734=\", \"75\u00d765=\"],\n  [\"49\u00d759=\", \"13\u00d737=\"],\n  [\"49\u00d742=\", \"75\u00d737=\"],\n  [\"98\u00d718=\", \"87\u00d752=\"],\n  [\"13\u00d734=\", \"45\u00d7100=\"],\n  [\"44\u00d710=\", \"38\u00d784=\"],\n  [\"28\u00d745=\", \"95\u00d793=\"],\n  [\"41\u00d711=\", \"96\u00d731=\"],\n  [\"26\u00d750=\", \"68\u00d726=\"],\n  [\"43\u00d740=\", \"77\u00d777=\"],\n  [\"96\u00d736=\", \"63\u00d798=\"],\n  [\"91\u00d721=\", \"59\u00d785=\"],\n  [\"15\u00d711=\", \"68\u00d766=\"],\n  [\"31\u00d756=\", \"17\u00d780=\"],\n  [\"100\u00d742=\", \"11\u00d781=\"],\n  [\"25\u00d757=\", \"83\u00d770=\"],\n  [\"94\u00d715=\", \"87\u00d740=\"],\n  [\"61\u00d746=\", \"74\u00d782=\"],\n  [\"100\u00d787=\", \"44\u00d758=\"],\n  [\"70\u00d782=\", \"26\u00d760=\"],\n  [\"51\u00d715=\", \"26\u00d791=\"],\n  [\"99\u00d765=\", \"23\u00d786=\"],\n  [\"21\u00d791=\", \"70\u00d786=\"],\n  [\"28\u00d713=\", \"41\u00d717=\"],\n  [\"86\u00d796=\", \"83\u00d736=\"],\n  [\"83\u00d781=\", \"19\u00d747=\"],\n  [\"74\u00d795=\", \"52\u00d757=\"],\n  [\"73\u00d785=\", \"76\u00d741=\"],\n  [\"52\u00d716=\", \"75\u00d748=\"],\n  [\"11\u00d773=\", \"14\u00d784=\"],\n  [\"57\u00d727=\", \"52\u00d764=\"],\n  [\"69\u00d760=\", \"71\u00d790=\"],\n  [\"91\u00d753=\", \"18\u00d781=\"],\n  [\"14\u00d711=\", \"50\u00d799=\"],\n  [\"62\u00d776=\", \"76\u00d760=\"],\n  [\"68\u00d755=\", \"66\u00d744=\"],\n  [\"51\u00d737=\", \"75\u00d746=\"],\n  [\"90\u00d788=\", \"95\u00d755=\"],\n  [\"50\u00d769=\", \"65\u00d780=\"],\n  [\"24\u00d721=\", \"49\u00d789=\"],\n  [\"37\u00d777=\", \"25\u00d775=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all 100 two-digit multiplication problems\n# in the practice-sheet table. Each (Old, New) pair is an exact, unique\n# value in this document, so Find/Execute with MatchCase + ReplaceAll is an\n# unambiguous single-hit replacement for every entry.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-04-26 Wednesday\", \"2023-04-27 Thursday\"),\n    @(\"96\u00d763=\", \"27\u00d755=\"),\n    @(\"88\u00d716=\", \"86\u00d789=\"),\n    @(\"46\u00d777=\", \"91\u00d744=\"),\n    @(\"87\u00d748=\", \"93\u00d742=\"),\n    @(\"50\u00d740=\", \"91\u00d757=\"),\n    @(\"84\u00d791=\", \"85\u00d797=\"),\n    @(\"48\u00d740=\", \"43\u00d7100=\"),\n    @(\"84\u00d745=\", \"12\u00d796=\"),\n    @(\"27\u00d794=\", \"22\u00d785=\"),\n    @(\"94\u00d747=\", \"25\u00d767=\"),\n    @(\"43\u00d777=\", \"99\u00d737=\"),\n    @(\"11\u00d745=\", \"48\u00d762=\"),\n    @(\"55\u00d751=\", \"87\u00d766=\"),\n    @(\"54\u00d717=\", \"71\u00d759=\"),\n    @(\"10\u00d785=\", \"11\u00d799=\"),\n    @(\"61\u00d796=\", \"10\u00d723=\"),\n    @(\"50\u00d754=\", \"30\u00d790=\"),\n    @(\"98\u00d711=\", \"60\u00d779=\"),\n    @(\"62\u00d746=\", \"55\u00d764=\"),\n    @(\"10\u00d726=\", \"74\u00d784=\"),\n    @(\"10\u00d729=\", \"39\u00d726=\"),\n    @(\"51\u00d745=\", \"39\u00d733=\"),\n    @(\"62\u00d733=\", \"33\u00d782=\"),\n    @(\"29\u00d733=\", \"57\u00d750=\"),\n    @(\"27\u00d722=\", \"13\u00d799=\"),\n    @(\"56\u00d758=\", \"48\u00d730=\"),\n    @(\"41\u00d760=\", \"11\u00d763=\"),\n    @(\"88\u00d736=\", \"83\u00d763=\"),\n    @(\"16\u00d753=\", \"20\u00d755=\"),\n    @(\"52\u00d767=\", \"96\u00d711=\"),\n    @(\"54\u00d764=\", \"94\u00d791=\"),\n    @(\"35\u00d771=\", \"74\u00d771=\"),\n    @(\"87\u00d783=\", \"54\u00d773=\"),\n    @(\"39\u00d720=\", \"49\u00d794=\"),\n    @(\"19\u00d748=\", \"36\u00d779=\"),\n    @(\"12\u00d749=\", \"47\u00d768=\"),\n    @(\"82\u00d758=\", \"97\u00d750=\"),\n    @(\"11\u00d754=\", \"74\u00d775=\"),\n    @(\"70\u00d711=\", \"59\u00d738=\"),\n    @(\"28\u00d714=\", \"71\u00d763=\"),\n    @(\"73\u00d734=\", \"65\u00d771=\"),\n    @(\"78\u00d762=\", \"52\u00d715=\"),\n    @(\"32\u00d745=\", \"73\u00d712=\"),\n    @(\"85\u00d727=\", \"80\u00d743=\"),\n    @(\"81\u00d793=\", \"95\u00d719=\"),\n    @(\"89\u00d740=\", \"96\u00d748=\"),\n    @(\"53\u00d744=\", \"12\u00d718=\"),\n    @(\"77\u00d723=\", \"20\u00d732=\"),\n    @(\"77\u00d715=\", \"42\u00d736=\"),\n    @(\"78\u00d782=\", \"61\u00d791=\"),\n    @(\"84\u00d769=\", \"38\u00d778=\"),\n    @(\"38\u00d796=\", \"39\u00d731=\"),\n    @(\"44\u00d767=\", \"60\u00d758=\"),\n    @(\"43\u00d788=\", \"67\u00d756=\"),\n    @(\"83\u00d798=\", \"37\u00d795=\"),\n    @(\"98\u00d784=\", \"64\u00d777=\"),\n    @(\"86\u00d774=\", \"71\u00d733=\"),\n    @(\"92\u00d794=\", \"34\u00d766=\"),\n    @(\"48\u00d720=\", \"30\u00d735=\"),\n    @(\"17\u00d734=\", \"75\u00d765=\"),\n    @(\"49\u00d759=\", \"13\u00d737=\"),\n    @(\"49\u00d742=\", \"75\u00d737=\"),\n    @(\"98\u00d718=\", \"87\u00d752=\"),\n    @(\"13\u00d734=\", \"45\u00d7100=\"),\n    @(\"44\u00d710=\", \"38\u00d784=\"),\n    @(\"28\u00d745=\", \"95\u00d793=\"),\n    @(\"41\u00d711=\", \"96\u00d731=\"),\n    @(\"26\u00d750=\", \"68\u00d726=\"),\n    @(\"43\u00d740=\", \"77\u00d777=\"),\n    @(\"96\u00d736=\", \"63\u00d798=\"),\n    @(\"91\u00d721=\", \"59\u00d785=\"),\n    @(\"15\u00d711=\", \"68\u00d766=\"),\n    @(\"31\u00d756=\", \"17\u00d780=\"),\n    @(\"100\u00d742=\", \"11\u00d781=\"),\n    @(\"25\u00d757=\", \"83\u00d770=\"),\n    @(\"94\u00d715=\", \"87\u00d740=\"),\n    @(\"61\u00d746=\", \"74\u00d782=\"),\n    @(\"100\u00d787=\", \"44\u00d758=\"),\n    @(\"70\u00d782=\", \"26\u00d760=\"),\n    @(\"51\u00d715=\", \"26\u00d791=\"),\n    @(\"99\u00d765=\", \"23\u00d786=\"),\n    @(\"21\u00d791=\", \"70\u00d786=\"),\n    @(\"28\u00d713=\", \"41\u00d717=\"),\n    @(\"86\u00d796=\", \"83\u00d736=\"),\n    @(\"83\u00d781=\", \"19\u00d747=\"),\n    @(\"74\u00d795=\", \"52\u00d757=\"),\n    @(\"73\u00d785=\", \"76\u00d741=\"),\n    @(\"52\u00d716=\", \"75\u00d748=\"),\n    @(\"11\u00d773=\", \"14\u00d784=\"),\n    @(\"57\u00d727=\", \"52\u00d764=\"),\n    @(\"69\u00d760=\", \"71\u00d790=\"),\n    @(\"91\u00d753=\", \"18\u00d781=\"),\n    @(\"14\u00d711=\", \"50\u00d799=\"),\n    @(\"62\u00d776=\", \"76\u00d760=\"),\n    @(\"68\u00d755=\", \"66\u00d744=\"),\n    @(\"51\u00d737=\", \"75\u00d746=\"),\n    @(\"90\u00d788=\", \"95\u00d755=\"),\n    @(\"50\u00d769=\", \"65\u00d780=\"),\n    @(\"24\u00d721=\", \"49\u00d789=\"),\n    @(\"37\u00d777=\", \"25\u00d775=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
